$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''65.106.20'
$ws.Range("E2").Value = '  +0.20%  '

# Row 3
$ws.Range("D3").Value = '''3.518.77'
$ws.Range("E3").Value = '  -1.33%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '''592.84'
$ws.Range("E5").Value = '  -1.25%  '

# Row 6
$ws.Range("D6").Value = '''134.30'
$ws.Range("E6").Value = '  -0.88%  '

# Row 7
$ws.Range("D7").Value = '''3.517.11'
$ws.Range("E7").Value = '  -1.31%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").Value = '''0.490'
$ws.Range("E9").Value = '  -0.91%  '

# Row 10
$ws.Range("D10").Value = '''0.125'
$ws.Range("E10").Value = '  +1.16%  '

# Row 11
$ws.Range("E11").Value = '  +3.04%  '

# Row 12
$ws.Range("D12").Value = '''0.387'
$ws.Range("E12").Value = '  +0.21%  '

# Row 13
$ws.Range("D13").Value = '''4.114.58'
$ws.Range("E13").Value = '  -1.42%  '

# Row 14
$ws.Range("D14").Value = '''27.76'
$ws.Range("E14").Value = '  +2.44%  '

# Row 15
$ws.Range("D15").Value = '''0.0000182'
$ws.Range("E15").Value = '  -0.52%  '

# Row 16
$ws.Range("E16").Value = '  +0.58%  '

# Row 17
$ws.Range("D17").Value = '''3.518.32'
$ws.Range("E17").Value = '  -1.69%  '

# Row 18
$ws.Range("D18").Value = '''65.056.89'
$ws.Range("E18").Value = '  -0.02%  '

# Row 19
$ws.Range("D19").Value = '''10.10'
$ws.Range("E19").Value = '  +0.22%  '

# Row 20
$ws.Range("D20").Value = '''14.34'
$ws.Range("E20").Value = '  -0.34%  '

# Row 21
$ws.Range("D21").Value = '''5.68'
$ws.Range("E21").Value = '  -3.14%  '

# Row 22
$ws.Range("D22").Value = '''392.49'
$ws.Range("E22").Value = '  +0.91%  '

# Row 23
$ws.Range("E23").Value = '  -0.20%  '

# Row 24
$ws.Range("D24").Value = '''3.659.76'
$ws.Range("E24").Value = '  -1.45%  '

# Row 25
$ws.Range("D25").Value = '''74.61'
$ws.Range("E25").Value = '  +0.58%  '

# Row 26
$ws.Range("E26").Value = '  +0.24%  '

# Row 27
$ws.Range("E27").Value = '  -4.25%  '

# Row 28
$ws.Range("D28").Value = '''1.62'
$ws.Range("E28").Value = '  +10.07%  '

# Row 29
$ws.Range("D29").Value = '''7.68'
$ws.Range("E29").Value = '  -0.90%  '

# Row 30
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.17%  '

# Row 31
$ws.Range("E31").Value = '  -1.09%  '

# Row 32
$ws.Range("D32").Value = '''8.31'
$ws.Range("E32").Value = '  -1.89%  '

# Row 33
$ws.Range("D33").Value = '''3.525.14'
$ws.Range("E33").Value = '  -1.11%  '

# Row 34
$ws.Range("D34").Value = '''24.11'
$ws.Range("E34").Value = '  +0.37%  '

# Row 35
$ws.Range("E35").Value = '  +0.00%  '

# Row 36
$ws.Range("E36").Value = '  +0.54%  '

# Row 37
$ws.Range("D37").Value = '''5.27'
$ws.Range("E37").Value = '  +5.18%  '

# Row 38
$ws.Range("E38").Value = '  +0.80%  '

# Row 39
$ws.Range("D39").Value = '''6.96'
$ws.Range("E39").Value = '  +0.16%  '

# Row 40
$ws.Range("D40").Value = '''168.31'
$ws.Range("E40").Value = '  -0.58%  '

# Row 41
$ws.Range("D41").Value = '''0.0807'
$ws.Range("E41").Value = '  -0.12%  '

# Row 42
$ws.Range("E42").Value = '  -0.59%  '

# Row 43
$ws.Range("E43").Value = '  +4.28%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''25.83'
$ws.Range("E44").Value = '  -5.27%  '

# Row 45
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '''42.95'
$ws.Range("E45").Value = '  +0.63%  '

# Row 46
$ws.Range("D46").Value = '''0.999'
$ws.Range("E46").Value = '  -0.03%  '

# Row 47
$ws.Range("D47").Value = '''4.44'
$ws.Range("E47").Value = '  -0.83%  '

# Row 48
$ws.Range("E48").Value = '  +0.84%  '

# Row 49
$ws.Range("E49").Value = '  -0.30%  '

# Row 50
$ws.Range("D50").Value = '''2.428.98'
$ws.Range("E50").Value = '  -2.65%  '

# Row 51
$ws.Range("D51").Value = '''0.906'
$ws.Range("E51").Value = '  +4.28%  '
